$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3651747703552246
$ws.Range("E2").Value = 228.0655785106464
$ws.Range("F2").Value = 0.007317468195793464
$ws.Range("G2").Value = 0.0060838686185918
$ws.Range("H2").Value = 0.005472297584576656
$ws.Range("I2").Value = 0.005472297584576656
$ws.Range("J2").Value = 0.005124048704155246
$ws.Range("K2").Value = 0.0050533156769138
$ws.Range("L2").Value = 0.004923794576204057
$ws.Range("M2").Value = 0.004866361063318834
$ws.Range("N2").Value = 0.00478863633447172
$ws.Range("O2").Value = 0.00478863633447172
$ws.Range("P2").Value = 0.004738758853910806
$ws.Range("Q2").Value = 0.004689256642850558
$ws.Range("R2").Value = 0.00463123908227971
$ws.Range("S2").Value = 0.00455334557268839
$ws.Range("T2").Value = 0.00455334557268839
$ws.Range("U2").Value = 0.004530252364020802
$ws.Range("V2").Value = 0.004500651061796064
$ws.Range("W2").Value = 0.004473453031385903
$ws.Range("X2").Value = 0.004466036502718127
$ws.Range("Y2").Value = 0.00444572277798531

$ws.Range("C3").Value = 0.4374756813049316
$ws.Range("E3").Value = 234.9109136191473
$ws.Range("F3").Value = 0.006854461259283044
$ws.Range("G3").Value = 0.005830092466304691
$ws.Range("H3").Value = 0.005548721327062377
$ws.Range("I3").Value = 0.005548721327062377
$ws.Range("J3").Value = 0.005548721327062377
$ws.Range("K3").Value = 0.005208458807797125
$ws.Range("L3").Value = 0.005139019219384589
$ws.Range("M3").Value = 0.005036010929215989
$ws.Range("N3").Value = 0.004991048391962181
$ws.Range("O3").Value = 0.004862559774808342
$ws.Range("P3").Value = 0.004862559774808342
$ws.Range("Q3").Value = 0.004768257003556728
$ws.Range("R3").Value = 0.004755565431364603
$ws.Range("S3").Value = 0.004685731512274595
$ws.Range("T3").Value = 0.004622698849021052
$ws.Range("U3").Value = 0.004622698849021052
$ws.Range("V3").Value = 0.004602385601989424
$ws.Range("W3").Value = 0.004598205479288067
$ws.Range("X3").Value = 0.004584981877047453
$ws.Range("Y3").Value = 0.004579160109535034

$ws.Range("C4").Value = 0.4062740802764893
$ws.Range("E4").Value = 226.2526246174166
$ws.Range("F4").Value = 0.007234092707530699
$ws.Range("G4").Value = 0.006284335820055877
$ws.Range("H4").Value = 0.005955803728992741
$ws.Range("I4").Value = 0.005473664674338815
$ws.Range("J4").Value = 0.005152333620179397
$ws.Range("K4").Value = 0.005114520711548476
$ws.Range("L4").Value = 0.004941327030793159
$ws.Range("M4").Value = 0.004732470209006531
$ws.Range("N4").Value = 0.004669679188246928
$ws.Range("O4").Value = 0.004669679188246928
$ws.Range("P4").Value = 0.004669679188246928
$ws.Range("Q4").Value = 0.004540147234083118
$ws.Range("R4").Value = 0.004522223036319399
$ws.Range("S4").Value = 0.00451472593612016
$ws.Range("T4").Value = 0.004510301587480662
$ws.Range("U4").Value = 0.004486600924080588
$ws.Range("V4").Value = 0.004470507402958979
$ws.Range("W4").Value = 0.004441119918072387
$ws.Range("X4").Value = 0.004422606320777598
$ws.Range("Y4").Value = 0.00441038254614847

$ws.Range("C5").Value = 0.3593747615814209
$ws.Range("E5").Value = 228.413687575845
$ws.Range("F5").Value = 0.007136695655527256
$ws.Range("G5").Value = 0.006254006836292634
$ws.Range("H5").Value = 0.00577683569309643
$ws.Range("I5").Value = 0.005578213332672863
$ws.Range("J5").Value = 0.00532610054434813
$ws.Range("K5").Value = 0.004926712895819371
$ws.Range("L5").Value = 0.004926712895819371
$ws.Range("M5").Value = 0.004798100816886849
$ws.Range("N5").Value = 0.004798100816886849
$ws.Range("O5").Value = 0.004748343004120983
$ws.Range("P5").Value = 0.004733619615442868
$ws.Range("Q5").Value = 0.004666523886144329
$ws.Range("R5").Value = 0.004598489349362841
$ws.Range("S5").Value = 0.004564065229739627
$ws.Range("T5").Value = 0.004535762942153607
$ws.Range("U5").Value = 0.004511304077250155
$ws.Range("V5").Value = 0.004473545028276416
$ws.Range("W5").Value = 0.004473545028276416
$ws.Range("X5").Value = 0.004462600901602556
$ws.Range("Y5").Value = 0.004452508529743565

$ws.Range("C6").Value = 0.3906266689300537
$ws.Range("E6").Value = 228.8258930396678
$ws.Range("F6").Value = 0.007337835762981425
$ws.Range("G6").Value = 0.006278761131290688
$ws.Range("H6").Value = 0.005701230951466495
$ws.Range("I6").Value = 0.00556853105381772
$ws.Range("J6").Value = 0.005439012448774522
$ws.Range("K6").Value = 0.00507026578790495
$ws.Range("L6").Value = 0.004980606111552063
$ws.Range("M6").Value = 0.004957508191375458
$ws.Range("N6").Value = 0.004957508191375458
$ws.Range("O6").Value = 0.004916640052115627
$ws.Range("P6").Value = 0.004819826142342576
$ws.Range("Q6").Value = 0.004764944392851003
$ws.Range("R6").Value = 0.004632902236873533
$ws.Range("S6").Value = 0.004621591858538886
$ws.Range("T6").Value = 0.004540205049836205
$ws.Range("U6").Value = 0.004540205049836205
$ws.Range("V6").Value = 0.004506427043875998
$ws.Range("W6").Value = 0.004484830475501378
$ws.Range("X6").Value = 0.004474916722642357
$ws.Range("Y6").Value = 0.004460543723970131

$ws.Range("C7").Value = 0.3749992847442627
$ws.Range("E7").Value = 229.8384626732277
$ws.Range("F7").Value = 0.007021633181949178
$ws.Range("G7").Value = 0.005991655727086207
$ws.Range("H7").Value = 0.005572918279442984
$ws.Range("I7").Value = 0.005572918279442984
$ws.Range("J7").Value = 0.005290486052376275
$ws.Range("K7").Value = 0.004999328347213399
$ws.Range("L7").Value = 0.004826088602703553
$ws.Range("M7").Value = 0.004726302550502507
$ws.Range("N7").Value = 0.004708407140275424
$ws.Range("O7").Value = 0.004618297977289237
$ws.Range("P7").Value = 0.004618297977289237
$ws.Range("Q7").Value = 0.004618297977289237
$ws.Range("R7").Value = 0.004616817058374321
$ws.Range("S7").Value = 0.004616817058374321
$ws.Range("T7").Value = 0.004580440321890731
$ws.Range("U7").Value = 0.004548930233975706
$ws.Range("V7").Value = 0.004543886409468342
$ws.Range("W7").Value = 0.004531915665950108
$ws.Range("X7").Value = 0.004497547734351
$ws.Range("Y7").Value = 0.00448028192345473

$ws.Range("C8").Value = 0.375
$ws.Range("E8").Value = 230.024404374024
$ws.Range("F8").Value = 0.007151321207335308
$ws.Range("G8").Value = 0.006083332396099871
$ws.Range("H8").Value = 0.005649262116197607
$ws.Range("I8").Value = 0.005649262116197607
$ws.Range("J8").Value = 0.005574198331899159
$ws.Range("K8").Value = 0.005259335052285191
$ws.Range("L8").Value = 0.00497184835210327
$ws.Range("M8").Value = 0.004819238123658451
$ws.Range("N8").Value = 0.004819238123658451
$ws.Range("O8").Value = 0.004775785842500259
$ws.Range("P8").Value = 0.004775785842500259
$ws.Range("Q8").Value = 0.004690306298701128
$ws.Range("R8").Value = 0.004628355234018071
$ws.Range("S8").Value = 0.004610931795815387
$ws.Range("T8").Value = 0.004593127558167885
$ws.Range("U8").Value = 0.004564216154279432
$ws.Range("V8").Value = 0.004526425456755759
$ws.Range("W8").Value = 0.004483906518012163
$ws.Range("X8").Value = 0.004483906518012163
$ws.Range("Y8").Value = 0.004483906518012163

$ws.Range("C9").Value = 0.3593747615814209
$ws.Range("E9").Value = 231.5373328050591
$ws.Range("F9").Value = 0.006877436701719059
$ws.Range("G9").Value = 0.005946869273461255
$ws.Range("H9").Value = 0.005520102474328858
$ws.Range("I9").Value = 0.005520102474328858
$ws.Range("J9").Value = 0.005421189613179872
$ws.Range("K9").Value = 0.005379677734621238
$ws.Range("L9").Value = 0.005181813291549269
$ws.Range("M9").Value = 0.004911443709806384
$ws.Range("N9").Value = 0.004911443709806384
$ws.Range("O9").Value = 0.004911443709806384
$ws.Range("P9").Value = 0.004854331496731922
$ws.Range("Q9").Value = 0.004833146608580668
$ws.Range("R9").Value = 0.00474748150220239
$ws.Range("S9").Value = 0.004688915594144442
$ws.Range("T9").Value = 0.004645436455360851
$ws.Range("U9").Value = 0.004584365466280129
$ws.Range("V9").Value = 0.004584365466280129
$ws.Range("W9").Value = 0.004567610292898167
$ws.Range("X9").Value = 0.004539308827067346
$ws.Range("Y9").Value = 0.004513398300293549

$ws.Range("C10").Value = 0.3593473434448242
$ws.Range("E10").Value = 226.9584485397936
$ws.Range("F10").Value = 0.007120334999264499
$ws.Range("G10").Value = 0.005877682751143919
$ws.Range("H10").Value = 0.005724512797741657
$ws.Range("I10").Value = 0.005207855372357735
$ws.Range("J10").Value = 0.005207855372357735
$ws.Range("K10").Value = 0.005207855372357735
$ws.Range("L10").Value = 0.005001061439661799
$ws.Range("M10").Value = 0.005001061439661799
$ws.Range("N10").Value = 0.004825634102335331
$ws.Range("O10").Value = 0.004825634102335331
$ws.Range("P10").Value = 0.004711551381330428
$ws.Range("Q10").Value = 0.004700925161609222
$ws.Range("R10").Value = 0.004678546245679674
$ws.Range("S10").Value = 0.004608545236211677
$ws.Range("T10").Value = 0.004551136719668111
$ws.Range("U10").Value = 0.004524417887490213
$ws.Range("V10").Value = 0.004507486716516283
$ws.Range("W10").Value = 0.004460200334910906
$ws.Range("X10").Value = 0.004448168102000749
$ws.Range("Y10").Value = 0.004424141297071999

$ws.Range("C11").Value = 0.3906278610229492
$ws.Range("E11").Value = 227.4362547215569
$ws.Range("F11").Value = 0.006948742911768419
$ws.Range("G11").Value = 0.006147391144375495
$ws.Range("H11").Value = 0.005640368498922989
$ws.Range("I11").Value = 0.005379255422924517
$ws.Range("J11").Value = 0.005048786077035681
$ws.Range("K11").Value = 0.004997053836085855
$ws.Range("L11").Value = 0.004997053836085855
$ws.Range("M11").Value = 0.004939224388412571
$ws.Range("N11").Value = 0.004886902177674242
$ws.Range("O11").Value = 0.004814849933728787
$ws.Range("P11").Value = 0.004703105578174377
$ws.Range("Q11").Value = 0.004703105578174377
$ws.Range("R11").Value = 0.004634195895024838
$ws.Range("S11").Value = 0.004620675524687608
$ws.Range("T11").Value = 0.004543923056693063
$ws.Range("U11").Value = 0.00451355129015704
$ws.Range("V11").Value = 0.004496860348228817
$ws.Range("W11").Value = 0.004476861988954355
$ws.Range("X11").Value = 0.004459350310820522
$ws.Range("Y11").Value = 0.004433455257730153
